# Update "想去人数" (number of people interested) figures across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3881
$ws1.Range("F11").Value = 1438
$ws1.Range("F13").Value = 2556

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 6

# Sheet "全部类型" (All types) - aggregate of the above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3881
$ws4.Range("F13").Value = 6
$ws4.Range("F14").Value = 1438
$ws4.Range("F16").Value = 2556
